$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the "ShibaInu" price text which contains a subscript-3 character (U+2083).
$sub3 = [char]0x2083
$shibaNew = "0.0" + "${sub3}" + "0752"

# Row-by-row updates: D = Price column, E = Volume(1h) column.
# D = $null means the Price cell is unchanged for that row (only Volume changed).
$updates = @(
    @{Row=2; D="26.618.97"; E="  -0.07%  "},
    @{Row=3; D="1.644.52"; E="  +0.72%  "},
    @{Row=4; D="1.00"; E="  +0.21%  "},
    @{Row=5; D="216.05"; E="  +1.38%  "},
    @{Row=6; D="0.503"; E="  +0.63%  "},
    @{Row=8; D="0.252"; E="  -0.14%  "},
    @{Row=9; D=$null; E="  +0.74%  "},
    @{Row=10; D=$null; E="  +0.60%  "},
    @{Row=11; D="0.0844"; E="  +0.05%  "},
    @{Row=12; D="1.874.77"; E="  +0.68%  "},
    @{Row=13; D=$null; E="  +3.27%  "},
    @{Row=14; D="1.646.65"; E="  +0.94%  "},
    @{Row=15; D="0.534"; E="  +1.84%  "},
    @{Row=16; D="66.13"; E="  +4.47%  "},
    @{Row=17; D="26.660.59"; E="  +0.09%  "},
    @{Row=18; D=$shibaNew; E="  +1.50%  "},
    @{Row=19; D="218.39"; E="  -0.48%  "},
    @{Row=20; D=$null; E="  +0.35%  "},
    @{Row=21; D=$null; E="  +2.11%  "},
    @{Row=22; D=$null; E="  +1.93%  "},
    @{Row=23; D="9.59"; E="  +1.56%  "},
    @{Row=24; D="2.11"; E="  +9.16%  "},
    @{Row=25; D="146.73"; E="  -1.36%  "},
    @{Row=26; D="1.00"; E="  +0.14%  "},
    @{Row=27; D=$null; E="  -0.19%  "},
    @{Row=28; D="7.14"; E="  +3.07%  "},
    @{Row=29; D="15.87"; E="  +2.30%  "},
    @{Row=30; D=$null; E="  +1.99%  "},
    @{Row=31; D=$null; E="  +0.86%  "},
    @{Row=32; D=$null; E="  +2.96%  "},
    @{Row=33; D=$null; E="  +2.48%  "},
    @{Row=34; D="1.279.21"; E="  +5.68%  "},
    @{Row=35; D=$null; E="  +2.14%  "},
    @{Row=36; D=$null; E="  +6.51%  "},
    @{Row=37; D=$null; E="  +0.32%  "},
    @{Row=38; D="0.528"; E="  +4.86%  "},
    @{Row=39; D=$null; E="  +1.94%  "},
    @{Row=40; D=$null; E="  +0.37%  "},
    @{Row=41; D=$null; E="  +2.02%  "},
    @{Row=42; D=$null; E="  -1.92%  "},
    @{Row=43; D=$null; E="  +0.71%  "},
    @{Row=44; D="1.786.11"; E="  +0.90%  "},
    @{Row=45; D="93.23"; E="  +0.24%  "},
    @{Row=46; D="59.72"; E="  +9.20%  "},
    @{Row=47; D=$null; E="  +4.15%  "},
    @{Row=48; D=$null; E="  +0.80%  "},
    @{Row=49; D=$null; E="  +1.96%  "},
    @{Row=50; D="0.0977"; E="  +3.67%  "},
    @{Row=51; D="0.408"; E="  -0.46%  "}
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($u.D -ne $null) {
        $dCell = $ws.Cells.Item($row, 4)
        # Force text storage so numeric-looking strings (e.g. "1.00") keep
        # their exact displayed form instead of being parsed into numbers,
        # then restore the cell's original style so no formatting changes
        # leak into the saved file.
        $origStyle = $dCell.Style
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
        $dCell.Style = $origStyle
    }

    $eCell = $ws.Cells.Item($row, 5)
    $eCell.Value = $u.E
}
